$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text number format on all changed cells to preserve exact string
# representation (e.g. trailing zeros, thousand-dot separators) as the
# source data stores these as plain text, not numeric values.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.500.30'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.393.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '564.15'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.23'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.539'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.397.87'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.56%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.35%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.847.12'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.528.02'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.35'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +15.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.374.68'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.90%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '325.37'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.88%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.80'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -8.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '64.46'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '557.45'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.26%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -11.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.522.53'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0911'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.30'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.30%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.129'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.31%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '153.52'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.24'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.07'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.34%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.57%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₆0280'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '142.92'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.30%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0500'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.97'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.61%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.30%  '
